$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the TNM override values for rows 3 (2030) and 4 (2035)
$ws.Range("C3").Value = 130000
$ws.Range("C4").Value = 160000

# Update the active selection cell as recorded in the sheet view
$ws.Range("I11").Select()
